$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture existing hyperlink targets (currently anchored at B3:E3) before
# removing them - the column insert below does not shift hyperlink anchors.
$hyperlinkTargets = @()
foreach ($hl in $ws.Hyperlinks) {
    $hyperlinkTargets += $hl.Address
}
$ws.Hyperlinks.Delete()

# Insert a new column before column A, shifting existing data (and formatting)
# one column to the right, so the hyperlinked cells become C3:F3.
$ws.Columns.Item(1).Insert()

# Recreate the hyperlinks on their new home cells, in the same order as before.
$destCells = @("C3", "D3", "E3", "F3")
for ($i = 0; $i -lt $hyperlinkTargets.Count; $i++) {
    $ws.Hyperlinks.Add($ws.Range($destCells[$i]), $hyperlinkTargets[$i]) | Out-Null
}
# Restore the original (non-duplicated) Hyperlink cell style on the cells.
$ws.Range("C3:F3").Style = "Hyperlink"

# Set header and data for the new column A
$ws.Range("A1").Value = "update"
$ws.Range("A2").Value = 20150809
$ws.Range("A3").Value = 20150809

# Update the selection to match the target (activeCell A4)
$ws.Range("A4").Select()
